$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# --- Row value/percentage updates ---
Set-TextValue $ws.Range("D2") "274.92"
Set-TextValue $ws.Range("E2") "-2.28%"
Set-TextValue $ws.Range("D3") "27.21"
Set-TextValue $ws.Range("E3") "1.44%"
Set-TextValue $ws.Range("D4") "4.754"
Set-TextValue $ws.Range("E4") "-3.79%"
Set-TextValue $ws.Range("D5") "0.06299"
Set-TextValue $ws.Range("E5") "-1.81%"
Set-TextValue $ws.Range("D6") "6.923"
Set-TextValue $ws.Range("E6") "-0.94%"
Set-TextValue $ws.Range("D7") "1.350"
Set-TextValue $ws.Range("E7") "34.74%"
Set-TextValue $ws.Range("D8") "0.8760"
Set-TextValue $ws.Range("E8") "-1.08%"
Set-TextValue $ws.Range("D9") "0.1507"
Set-TextValue $ws.Range("E9") "1.32%"
Set-TextValue $ws.Range("D10") "0.05026"
Set-TextValue $ws.Range("E10") "-3.20%"
Set-TextValue $ws.Range("D11") "0.07580"
Set-TextValue $ws.Range("E11") "2.37%"
Set-TextValue $ws.Range("D12") "0.02971"
Set-TextValue $ws.Range("E12") "-4.34%"
Set-TextValue $ws.Range("D13") "0.09000"
Set-TextValue $ws.Range("E13") "-0.50%"
Set-TextValue $ws.Range("D14") "0.001570"
Set-TextValue $ws.Range("E14") "-0.46%"
Set-TextValue $ws.Range("D15") "0.0006357"
Set-TextValue $ws.Range("E15") "0.67%"
Set-TextValue $ws.Range("D16") "0.006026"
Set-TextValue $ws.Range("E16") "0.23%"
Set-TextValue $ws.Range("E17") "-1.37%"
Set-TextValue $ws.Range("D18") "3.296"
Set-TextValue $ws.Range("E18") "-1.69%"
Set-TextValue $ws.Range("E20") "0.12%"
Set-TextValue $ws.Range("D21") "0.1344"
Set-TextValue $ws.Range("E21") "1.00%"
Set-TextValue $ws.Range("D22") "3.923"
Set-TextValue $ws.Range("E22") "-0.17%"
Set-TextValue $ws.Range("E23") "1.02%"
Set-TextValue $ws.Range("E24") "-0.52%"
Set-TextValue $ws.Range("D25") "0.003839"
Set-TextValue $ws.Range("E25") "4.01%"
Set-TextValue $ws.Range("E26") "0.07%"
Set-TextValue $ws.Range("D27") "0.0001936"
Set-TextValue $ws.Range("E27") "14.33%"
Set-TextValue $ws.Range("E40") "0.04%"
Set-TextValue $ws.Range("D41") "0.006844"
Set-TextValue $ws.Range("E41") "2.93%"
Set-TextValue $ws.Range("D42") "0.1172"
Set-TextValue $ws.Range("E42") "-0.62%"
Set-TextValue $ws.Range("D43") "0.002070"
Set-TextValue $ws.Range("E43") "-12.22%"
Set-TextValue $ws.Range("D44") "0.01151"
Set-TextValue $ws.Range("E44") "-8.26%"
Set-TextValue $ws.Range("D45") "0.00005159"
Set-TextValue $ws.Range("E45") "-1.86%"

# --- Row 46/47: BOLO and CoinbaseStockToken swap places ---
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws.Range("D46") "0.02299"
Set-TextValue $ws.Range("E46") "2.24%"

$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws.Range("D47") "1.490"
Set-TextValue $ws.Range("E47") "-36.76%"
